$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH248"
$ws.Range("C2").Value = "THE SPREAD OF CHGRISTIANITY IN AFRICA DURING THE 19TH AND 20TH CENTURIES"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"

$rng = $ws.Range("A2:H2")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 10
$rng.Font.ThemeColor = 1
